$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row above row 13 (this shifts existing rows 13.. down by one)
$ws.Rows.Item(13).Insert()

# Give the new row the same formatting as the (now shifted) row below it, which
# carries the same "plain todo row" format (border etc.) that the inserted row
# should have. Using copy/paste-formats reuses the existing cell style instead
# of fabricating a new one.
$ws.Range("A14").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the new cell's value (goes into shared strings as a new unique string)
$ws.Range("A13").Value = "Themeneingabe Planung/Deliverables/etc. anpassen??"

# Restore the selection to A13 as recorded in the saved workbook
$ws.Range("A13").Select()
